# Apply updated dSF (column F) values for the maton_phil 2024 save data.
# This corresponds to a "repull data, push all data, mean calculation" update
# where several rows' final swing displacement (dSF) values were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = -5
    4  = 3
    7  = -2
    15 = -2
    24 = -1
    27 = -11
    28 = -4
    31 = -6
    37 = 1
    43 = 3
    50 = -7
    51 = -5
    53 = -2
    59 = 1
    71 = -4
    75 = 2
    76 = -6
    80 = 0
    81 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
